# FrameworkConstants test data: row 3's scenario label referenced the wrong
# destination ("Mayiladuthurai") even though the row's destination city is
# Karaikal. Fix the scenario label and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the mismatched scenario name for row 3 (destination is Karaikal).
$ws.Range("A3").Value = "From Chennai to Karaikal"

# Update the active cell/selection on the sheet.
$ws.Range("C14").Select()
